$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.1035853774462193
$ws.Range("E2").Value = 0.002950897072710124

$ws.Range("D3").Value = 0.1047916968031245
$ws.Range("E3").Value = -0.002499519323207244

$ws.Range("D4").Value = 0.1172997044905044
$ws.Range("E4").Value = -0.0002510460251046176

$ws.Range("D5").Value = 0.1376295867688895
$ws.Range("E5").Value = 0.00604101080575159

$ws.Range("D6").Value = 0.1350233858913159
$ws.Range("E6").Value = 0.002346729246113233

$ws.Range("D7").Value = 0.1451612903225806
$ws.Range("E7").Value = 0.007972097658196509

$ws.Range("D8").Value = 0.1277355293339671
$ws.Range("E8").Value = 0.007751937984496138

$ws.Range("D9").Value = 0.1287734289433986
$ws.Range("E9").Value = 0.009982989581118318

$ws.Range("E10").Value = 0.004595560124816434
